$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Translatable_Site_labels")

# Insert a new row at position 19 (shifts existing rows 19+ down by one)
$ws.Rows.Item(19).Insert()

# Populate the new row with the "No known works" / worksUnknown label
# (order matters for shared-string table placement: "No known works" must
# land before "worksUnknown" to match the original authoring order)
$ws.Range("C19").Value = "No known works"
$ws.Range("B19").Value = "worksUnknown"
$ws.Range("D19").Value = "Author"
$ws.Range("E19").Formula = "=_xlfn.CONCAT("""",B19,"" : '"",C19,""',"")"

# Row 10's E cell gets overwritten with a literal value
$ws.Range("E10").Value = "Admin"

# Leave the sheet selection on the full formula column, matching the
# author's final selection state
$ws.Activate() | Out-Null
$ws.Range("E1:E64").Select() | Out-Null
